$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 12: label + two formulas (angle of inclination calc for 55 Cnc e)
$ws.Range("B12").Value = "i="
$ws.Range("C12").Formula = "=ACOS(((C5*696340000)+(1.6*69911000))/(D2*149597870700))"
$ws.Range("D12").Formula = "=C12*180/PI()"

# Update the active selection to match the saved state in the diff
$ws.Range("G21").Select()
